# Update "想去人数" (want-to-go count) figures on the sheets that carry
# exhibition data: "展览" and "全部类型" both list the same four rows
# (F2, F5, F6, F8) that need to be bumped by a small amount.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2" = 2284
    "F5" = 1095
    "F6" = 865
    "F8" = 5853
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
